# Update the "grid_cell" column (AG, rows 4-26) of the wind ~FI_T table
# on the "solar" worksheet - each wind process row is being re-paired with
# a different CHE grid cell label than it had before.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$values = @(
    "CHE_20",
    "CHE_1",
    "CHE_6",
    "CHE_11",
    "CHE_15",
    "CHE_25",
    "CHE_13",
    "CHE_0",
    "CHE_14",
    "CHE_18",
    "CHE_3",
    "CHE_24",
    "CHE_5",
    "CHE_8",
    "CHE_7",
    "CHE_10",
    "CHE_22",
    "CHE_17",
    "CHE_19",
    "CHE_12",
    "CHE_21",
    "CHE_9",
    "CHE_4"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 4 + $i
    $ws.Range("AG" + $row).Value = $values[$i]
}
